# Load 3: Load TA File + ML (Line Regression)
# Updates row 18 (last data row) with refreshed metrics values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = 22.83625233803775
$ws.Range("C18").Value = 82
$ws.Range("D18").Value = 73.17073170731707
$ws.Range("E18").Value = -1099.879599999999
$ws.Range("F18").Value = -2.40343866756536
$ws.Range("G18").Value = 380289.5273999991
$ws.Range("H18").Value = 831.0023706443095
$ws.Range("I18").Value = 34.0305737928054
$ws.Range("J18").Value = -20.73170731707317
